$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.657.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.05%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.487.88'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '546.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.73%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.68%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.486.75'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.39%  '

$ws.Range("E10").Value = '  -10.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.153'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.34'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -8.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -7.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.933.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.549.81'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.88%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000162'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -9.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.475.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.84%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '318.76'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.98%  '

$ws.Range("E26").Value = '  -5.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.606.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.48%  '

$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.17%  '

$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '537.37'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.47'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.71%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.72%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.54'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.146'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.57%  '

$ws.Range("E34").Value = '  -8.33%  '

$ws.Range("E35").Value = '  -9.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.76'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -11.53%  '

$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.12%  '

$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.84'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.374'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.40'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '144.12'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.51'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '148.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.18%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.89'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -9.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0531'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.64%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.587'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0938'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.35%  '
